$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = "calleeMd1.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-23 15:11:12"

$wsOverview.Range("A3").Value = "calleeMd2.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-23 15:11:12"

$wsOverview.Range("A4").Value = "callerMd1.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-03-23 15:11:12"

$wsOverview.Range("A5").Value = "callerMd2.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-23 15:11:12"

$e2eBase = "https://github.com/OpenLocalizationTest/oltest/blob/c6f2f0dfea0910a7c041a57465224fa03992c0ba/e2e/"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($e2eBase + "calleeMd1.md"), [Type]::Missing, [Type]::Missing, "calleeMd1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($e2eBase + "calleeMd2.md"), [Type]::Missing, [Type]::Missing, "calleeMd2.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), ($e2eBase + "callerMd1.md"), [Type]::Missing, [Type]::Missing, "callerMd1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), ($e2eBase + "callerMd2.md"), [Type]::Missing, [Type]::Missing, "callerMd2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

$zhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11c131b51921843288c529170e6383a47c6d75f0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"

# Row 2 - calleeMd1.md
$wsZh.Range("A2").Value = "calleeMd1.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-23 15:11:07"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("J2").Value = "Include"
$wsZh.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

# Row 3 - calleeMd2.md
$wsZh.Range("A3").Value = "calleeMd2.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-23 15:11:07"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("J3").Value = "Include"
$wsZh.Range("K3").Value = "e2e\callerMd1.md"

# Row 4 - callerMd1.md
$wsZh.Range("A4").Value = "callerMd1.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"
$wsZh.Range("E4").Value = "2016-03-23 15:11:07"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$wsZh.Range("J4").Value = "Include"

# Row 5 - callerMd2.md (new row)
$wsZh.Range("A5").Value = "callerMd2.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"
$wsZh.Range("E5").Value = "2016-03-23 15:11:07"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value = "e2e\calleeMd1.md"
$wsZh.Range("J5").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($e2eBase + "calleeMd1.md"), [Type]::Missing, [Type]::Missing, "calleeMd1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), ($zhBase + "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($e2eBase + "calleeMd2.md"), [Type]::Missing, [Type]::Missing, "calleeMd2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), ($zhBase + "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ($e2eBase + "callerMd1.md"), [Type]::Missing, [Type]::Missing, "callerMd1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), ($zhBase + "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ($e2eBase + "callerMd2.md"), [Type]::Missing, [Type]::Missing, "callerMd2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), ($zhBase + "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

$deBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2e9be323aa7bccaa2ff10641158fb1b16d5576b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# Row 2 - calleeMd1.md
$wsDe.Range("A2").Value = "calleeMd1.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-23 15:11:12"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("J2").Value = "Include"
$wsDe.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

# Row 3 - calleeMd2.md
$wsDe.Range("A3").Value = "calleeMd2.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-23 15:11:12"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("J3").Value = "Include"
$wsDe.Range("K3").Value = "e2e\callerMd1.md"

# Row 4 - callerMd1.md
$wsDe.Range("A4").Value = "callerMd1.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"
$wsDe.Range("E4").Value = "2016-03-23 15:11:12"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$wsDe.Range("J4").Value = "Include"

# Row 5 - callerMd2.md (new row)
$wsDe.Range("A5").Value = "callerMd2.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"
$wsDe.Range("E5").Value = "2016-03-23 15:11:12"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value = "e2e\calleeMd1.md"
$wsDe.Range("J5").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($e2eBase + "calleeMd1.md"), [Type]::Missing, [Type]::Missing, "calleeMd1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), ($deBase + "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"), [Type]::Missing, [Type]::Missing, "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($e2eBase + "calleeMd2.md"), [Type]::Missing, [Type]::Missing, "calleeMd2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), ($deBase + "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"), [Type]::Missing, [Type]::Missing, "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ($e2eBase + "callerMd1.md"), [Type]::Missing, [Type]::Missing, "callerMd1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), ($deBase + "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"), [Type]::Missing, [Type]::Missing, "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ($e2eBase + "callerMd2.md"), [Type]::Missing, [Type]::Missing, "callerMd2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), ($deBase + "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"), [Type]::Missing, [Type]::Missing, "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf") | Out-Null
